$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3634.1428
$ws.Range("I64").Value = 2957.3333
$ws.Range("J64").Value = 4141.75
$ws.Range("K64").Value = 2957.3333
$ws.Range("L64").Value = 4141.75
$ws.Range("M64").Value = -2709.3333
$ws.Range("N64").Value = -4637.75

$ws.Range("H67").Value = 3634.1428
$ws.Range("I67").Value = 2957.3333
$ws.Range("J67").Value = 4141.75
$ws.Range("K67").Value = 2957.3333
$ws.Range("L67").Value = 4141.75
$ws.Range("M67").Value = -2099.3333
$ws.Range("N67").Value = -5857.75

$ws.Range("H76").Value = 2649163.8
$ws.Range("I76").Value = 2852518.2
$ws.Range("J76").Value = 5555
$ws.Range("K76").Value = 2852518.2
$ws.Range("L76").Value = 5555
$ws.Range("M76").Value = -2852203.2
$ws.Range("N76").Value = -6185

$ws.Range("H79").Value = 2649163.8
$ws.Range("I79").Value = 2852518.2
$ws.Range("J79").Value = 5555
$ws.Range("K79").Value = 2852518.2
$ws.Range("L79").Value = 5555
$ws.Range("M79").Value = -2851426.2
$ws.Range("N79").Value = -7739

$ws.Range("H135").Value = 2412
$ws.Range("I135").Value = 1327.4138
$ws.Range("J135").Value = 6343.625
$ws.Range("K135").Value = 11946.7242
$ws.Range("L135").Value = 57092.625
$ws.Range("M135").Value = -9411.724200000001
$ws.Range("N135").Value = -62162.625

$ws.Range("H137").Value = 1009.7143
$ws.Range("J137").Value = 1646.0769
$ws.Range("L137").Value = 4938.2307
$ws.Range("N137").Value = -10038.2307

$ws.Range("H139").Value = 70416
$ws.Range("J139").Value = 70416
$ws.Range("L139").Value = 70416
$ws.Range("N139").Value = -80696

$ws.Range("H140").Value = 71322.73
$ws.Range("J140").Value = 90568.75
$ws.Range("L140").Value = 90568.75
$ws.Range("N140").Value = -100928.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 500071.12
$ws.Range("I32").Value = 7004.0205
$ws.Range("J32").Value = 3184547.5
$ws.Range("K32").Value = 7004.0205
$ws.Range("L32").Value = 3184547.5
$ws.Range("M32").Value = -6717.0205
$ws.Range("N32").Value = -3185121.5

$ws.Range("H74").Value = 1017.375
$ws.Range("I74").Value = 697.55554
$ws.Range("J74").Value = 1428.5714
$ws.Range("K74").Value = 697.55554
$ws.Range("L74").Value = 1428.5714
$ws.Range("M74").Value = 176.44446
$ws.Range("N74").Value = -3176.5714

$ws.Range("H77").Value = 1017.375
$ws.Range("I77").Value = 697.55554
$ws.Range("J77").Value = 1428.5714
$ws.Range("K77").Value = 3487.7777
$ws.Range("L77").Value = 7142.857
$ws.Range("M77").Value = 880.2223000000004
$ws.Range("N77").Value = -15878.857

$ws.Range("H139").Value = 60966.668
$ws.Range("J139").Value = 60966.668
$ws.Range("L139").Value = 60966.668
$ws.Range("N139").Value = -71246.66800000001

$ws.Range("H140").Value = 94454.55
$ws.Range("J140").Value = 94454.55
$ws.Range("L140").Value = 94454.55
$ws.Range("N140").Value = -104814.55

$ws.Range("H141").Value = 59793.332
$ws.Range("J141").Value = 62492.855
$ws.Range("L141").Value = 62492.855
$ws.Range("N141").Value = -72852.85500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1753.9333
$ws.Range("I105").Value = 1775.75
$ws.Range("J105").Value = 1666.6666
$ws.Range("K105").Value = 1775.75
$ws.Range("L105").Value = 1666.6666
$ws.Range("M105").Value = -28.75
$ws.Range("N105").Value = -5160.6666

$ws.Range("H140").Value = 59200
$ws.Range("J140").Value = 59200
$ws.Range("L140").Value = 59200
$ws.Range("N140").Value = -69560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 5186.6665
$ws.Range("I23").Value = 5250
$ws.Range("J23").Value = 5168.5713
$ws.Range("K23").Value = 5250
$ws.Range("L23").Value = 5168.5713
$ws.Range("M23").Value = -5010
$ws.Range("N23").Value = -5648.5713

$ws.Range("H27").Value = 5186.6665
$ws.Range("I27").Value = 5250
$ws.Range("J27").Value = 5168.5713
$ws.Range("K27").Value = 5250
$ws.Range("L27").Value = 5168.5713
$ws.Range("M27").Value = -5058
$ws.Range("N27").Value = -5552.5713

$ws.Range("H31").Value = 8381.798000000001
$ws.Range("I31").Value = 2865.558
$ws.Range("J31").Value = 13032.745
$ws.Range("K31").Value = 2865.558
$ws.Range("L31").Value = 13032.745
$ws.Range("M31").Value = -2570.558
$ws.Range("N31").Value = -13622.745

$ws.Range("H34").Value = 8381.798000000001
$ws.Range("I34").Value = 2865.558
$ws.Range("J34").Value = 13032.745
$ws.Range("K34").Value = 2865.558
$ws.Range("L34").Value = 13032.745
$ws.Range("M34").Value = -2663.558
$ws.Range("N34").Value = -13436.745

$ws.Range("H62").Value = 3896.1765
$ws.Range("I62").Value = 3864.2307
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3864.2307
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3240.2307
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 3896.1765
$ws.Range("I65").Value = 3864.2307
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 19321.1535
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -16201.1535
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14809877
$ws.Range("I70").Value = 25574996
$ws.Range("J70").Value = 7838.625
$ws.Range("K70").Value = 25574996
$ws.Range("L70").Value = 7838.625
$ws.Range("M70").Value = -25574726
$ws.Range("N70").Value = -8378.625

$ws.Range("H73").Value = 14809877
$ws.Range("I73").Value = 25574996
$ws.Range("J73").Value = 7838.625
$ws.Range("K73").Value = 25574996
$ws.Range("L73").Value = 7838.625
$ws.Range("M73").Value = -25574060
$ws.Range("N73").Value = -9710.625

$ws.Range("H80").Value = 72257
$ws.Range("I80").Value = 103573.82
$ws.Range("J80").Value = 3360
$ws.Range("K80").Value = 103573.82
$ws.Range("L80").Value = 3360
$ws.Range("M80").Value = -102575.82
$ws.Range("N80").Value = -5356

$ws.Range("H83").Value = 72257
$ws.Range("I83").Value = 103573.82
$ws.Range("J83").Value = 3360
$ws.Range("K83").Value = 517869.1
$ws.Range("L83").Value = 16800
$ws.Range("M83").Value = -512877.1
$ws.Range("N83").Value = -26784

$ws.Range("H138").Value = 68550
$ws.Range("J138").Value = 68550
$ws.Range("L138").Value = 68550
$ws.Range("N138").Value = -78830

$ws.Range("H140").Value = 89979
$ws.Range("J140").Value = 89979
$ws.Range("L140").Value = 89979
$ws.Range("N140").Value = -100339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2787.8386
$ws.Range("I7").Value = 2330.5217
$ws.Range("J7").Value = 4102.625
$ws.Range("K7").Value = 2330.5217
$ws.Range("L7").Value = 4102.625
$ws.Range("M7").Value = -2218.5217
$ws.Range("N7").Value = -4326.625

$ws.Range("H40").Value = 2459.2222
$ws.Range("I40").Value = 2227.2666
$ws.Range("J40").Value = 3619
$ws.Range("K40").Value = 2227.2666
$ws.Range("L40").Value = 3619
$ws.Range("M40").Value = -2091.2666
$ws.Range("N40").Value = -3891

$ws.Range("H61").Value = 2366.7812
$ws.Range("I61").Value = 1832.174
$ws.Range("J61").Value = 3733
$ws.Range("K61").Value = 1832.174
$ws.Range("L61").Value = 3733
$ws.Range("M61").Value = -1630.174
$ws.Range("N61").Value = -4137

$ws.Range("H113").Value = 2366.7812
$ws.Range("I113").Value = 1832.174
$ws.Range("J113").Value = 3733
$ws.Range("K113").Value = 1832.174
$ws.Range("L113").Value = 3733
$ws.Range("M113").Value = 337.826
$ws.Range("N113").Value = -8073

$ws.Range("H126").Value = 2787.8386
$ws.Range("I126").Value = 2330.5217
$ws.Range("J126").Value = 4102.625
$ws.Range("K126").Value = 6991.5651
$ws.Range("L126").Value = 12307.875
$ws.Range("M126").Value = -4521.5651
$ws.Range("N126").Value = -17247.875

$ws.Range("H138").Value = 52252.77
$ws.Range("J138").Value = 52252.77
$ws.Range("L138").Value = 52252.77
$ws.Range("N138").Value = -62532.77

$ws.Range("H140").Value = 58920
$ws.Range("J140").Value = 59400
$ws.Range("L140").Value = 59400
$ws.Range("N140").Value = -69760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 604.2646999999999
$ws.Range("I107").Value = 615.3913
$ws.Range("J107").Value = 581
$ws.Range("K107").Value = 1846.1739
$ws.Range("L107").Value = 1743
$ws.Range("M107").Value = 73.8261
$ws.Range("N107").Value = -5583

$ws.Range("H122").Value = 2894.4736
$ws.Range("I122").Value = 1935
$ws.Range("J122").Value = 5581
$ws.Range("K122").Value = 5805
$ws.Range("L122").Value = 16743
$ws.Range("M122").Value = -3355
$ws.Range("N122").Value = -21643

$ws.Range("H139").Value = 69900
$ws.Range("J139").Value = 69900
$ws.Range("L139").Value = 69900
$ws.Range("N139").Value = -80180

$ws.Range("H140").Value = 29305.8
$ws.Range("J140").Value = 29305.8
$ws.Range("L140").Value = 29305.8
$ws.Range("N140").Value = -39665.8
